$d = $word.ActiveDocument

# Locate the "Notes of the meeting:" paragraph - the new content is inserted
# right after it.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs($i).Range.Text.TrimEnd("`r")
    if ($text -eq "Notes of the meeting:") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find the 'Notes of the meeting:' paragraph"
}

$anchor = $d.Paragraphs($anchorIndex)

# Insert a new empty paragraph right after the anchor, then fill it in with
# the first new line of text.
$anchor.Range.InsertParagraphAfter()
$line1 = $d.Paragraphs($anchorIndex + 1)
$line1.Range.Text = "Heroku look at this for the database (digital ocean already works so forget it)"

# Insert the second new paragraph (with text) after the first one.
$line1.Range.InsertParagraphAfter()
$line2 = $d.Paragraphs($anchorIndex + 2)
$line2.Range.Text = "Unit tests"

# Insert the third new paragraph, which stays empty (matching the blank
# paragraph introduced in the diff).
$line2.Range.InsertParagraphAfter()
